$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D11").Value = "A"
$ws.Range("D19").Value = "A"
$ws.Range("D30").Value = "A"
$ws.Range("D38").Value = "BBB"
$ws.Range("D46").Value = "B"
$ws.Range("D58").Value = "BB"
$ws.Range("D59").Value = "B"
$ws.Range("D69").Value = "BB"
$ws.Range("D70").Value = "BB"
$ws.Range("D74").Value = "B"
$ws.Range("D80").Value = "B"
$ws.Range("D81").Value = "AAA"
$ws.Range("D87").Value = "BBB"
$ws.Range("D94").Value = "A"
$ws.Range("D99").Value = "A"
$ws.Range("D116").Value = "BBB"
$ws.Range("D124").Value = "BBB"
$ws.Range("D136").Value = "BBB"
$ws.Range("D143").Value = "BBB"
$ws.Range("D146").Value = "BB"
$ws.Range("D159").Value = "A"
$ws.Range("D160").Value = "A"
$ws.Range("D166").Value = "A"
$ws.Range("D167").Value = "A"
$ws.Range("D201").Value = "BBB"
$ws.Range("D204").Value = "A"
$ws.Range("D207").Value = "BBB"
$ws.Range("D233").Value = "BB"
$ws.Range("D281").Value = "A"
$ws.Range("D286").Value = "A"
$ws.Range("D287").Value = "A"
$ws.Range("D291").Value = "BBB"
$ws.Range("D292").Value = "BBB"
$ws.Range("D293").Value = "BBB"
$ws.Range("D302").Value = "AA"
$ws.Range("D306").Value = "A"
$ws.Range("D307").Value = "BBB"
$ws.Range("D308").Value = "BBB"
$ws.Range("D315").Value = "BB"
$ws.Range("D317").Value = "BB"
$ws.Range("D323").Value = "B"
$ws.Range("D330").Value = "A"
$ws.Range("D342").Value = "BBB"
$ws.Range("D348").Value = "BBB"
$ws.Range("D349").Value = "BB"
$ws.Range("D355").Value = "BB"
$ws.Range("D362").Value = "B"
$ws.Range("D375").Value = "A"
$ws.Range("D383").Value = "B"
$ws.Range("D387").Value = "BBB"
$ws.Range("D393").Value = "A"
